# Applique les changements décrits par le diff :
#  - Ajoute une nouvelle ligne de dépense "Condensateurs" (ligne 16)
#  - Met à jour la formule de L7 pour inclure H16
#  - Change la cellule sélectionnée de M14 à L8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvelle ligne d'article : Condensateurs, quantité 1, prix unitaire 9$, achat par MS
$ws.Range("C16").Value = "Condensateurs"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 9
$ws.Range("G16").Value = "MS"
$ws.Range("J16").Formula = "=F16"

# Mise à jour de la formule du total MS (L7) pour inclure la nouvelle ligne (H16)
$ws.Range("L7").Formula = "=H4+H5+H10+H13+H14+H16"

# Mise à jour de la cellule sélectionnée dans la vue
$ws.Range("L8").Select()
